$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing six
# header columns (A:F) one place to the right (B:G) and the old
# column-width formatting shifts with them.
$ws.Columns("A:A").Insert()

# New first two header columns.
$ws.Range("A1").Value = "Pipeline Component"
$ws.Range("B1").Value = "Workflow / CLT?"

# A new data row describing the 01-qc-se.cwl pipeline component (row 2
# is intentionally left blank, matching the target layout).
$ws.Range("A3").Value = "01-qc-se.cwl"
$ws.Range("B3").Value = "Workflow"
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = "No"
$ws.Range("E3").Value = "GGR"
$ws.Range("F3").Value = "No"
$ws.Range("G3").Value = "No"

# Re-set explicit column widths to the new layout.
$ws.Columns("A:B").ColumnWidth = 24.5
$ws.Columns("C:C").ColumnWidth = 12.833333333333334
$ws.Columns("D:D").ColumnWidth = 27
$ws.Columns("E:E").ColumnWidth = 29.166666666666668
$ws.Columns("F:F").ColumnWidth = 15.5
$ws.Columns("G:G").ColumnWidth = 16.666666666666668

# Header row: bigger bold font and a taller row to match.
$ws.Range("A1:G1").Font.Size = 14
$ws.Rows(1).RowHeight = 19

# Move the active selection the way the author left it.
$ws.Range("A4").Select() | Out-Null
